# Update the "Förändrad" (changed) date column for rows 2-12 from 2023-10-25 (45224)
# to 2023-11-03 (45233), matching the automatic data refresh captured in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 12; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45224) {
        $cell.Value2 = 45233
    }
}
